$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '24.883.11'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  +0.01%  '
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.707.09'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  -0.09%  '
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = '1.0000'
$dCell.Style = "Normal"
$ws.Range("E4").Value = '  -0.48%  '
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '317.92'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '1.000'
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  -0.25%  '
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '0.4004'
$dCell.Style = "Normal"
$ws.Range("E7").Value = '  +0.90%  '
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '0.4085'
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  -0.49%  '
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '1.487'
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  -1.72%  '
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  -0.35%  '
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '53.64'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  +0.15%  '
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '0.08856'
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  -1.59%  '
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '26.36'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  +8.14%  '
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '7.507'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  -2.58%  '
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '8.156'
$dCell.Style = "Normal"
$ws.Range("E15").Value = '  -0.78%  '
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '0.00001362'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  -1.77%  '
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '1.706.26'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  -1.20%  '
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '97.10'
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  -3.10%  '
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '0.07187'
$dCell.Style = "Normal"
$ws.Range("E19").Value = '  +0.27%  '
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '21.31'
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  +6.01%  '
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '7.302'
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  -3.40%  '
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9980'
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  -1.12%  '
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '14.42'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  -1.05%  '
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '24.869.34'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '2.951'
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  -4.59%  '
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '2.326'
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  -0.81%  '
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '23.34'
$dCell.Style = "Normal"
$ws.Range("E27").Value = '  +1.14%  '
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '6.312'
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  +20.54%  '
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '166.97'
$dCell.Style = "Normal"
$ws.Range("E29").Value = '  +0.12%  '
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '146.47'
$dCell.Style = "Normal"
$ws.Range("E30").Value = '  +5.09%  '
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '8.427'
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  -9.61%  '
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '2.238'
$dCell.Style = "Normal"
$ws.Range("E32").Value = '  +13.78%  '
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '1.893.80'
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  -1.03%  '
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '0.08858'
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  -3.13%  '
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '0.03208'
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  +4.91%  '
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '7.215'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  -8.73%  '
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '1.042'
$dCell.Style = "Normal"
$ws.Range("E37").Value = '  -4.11%  '
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.2891'
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  +2.46%  '
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '0.8570'
$dCell.Style = "Normal"
$ws.Range("E39").Value = '  +8.65%  '
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '10.90'
$dCell.Style = "Normal"
$ws.Range("E40").Value = '  -2.25%  '
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '0.09359'
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  +0.20%  '
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '14.24'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  -2.59%  '
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '1.471'
$dCell.Style = "Normal"
$ws.Range("E43").Value = '  -1.11%  '
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '17.54'
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  +4.66%  '
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '2.726'
$dCell.Style = "Normal"
$ws.Range("E45").Value = '  +2.64%  '
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '0.7487'
$dCell.Style = "Normal"
$ws.Range("E46").Value = '  +2.07%  '
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '4.245'
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  -0.76%  '
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '1.405'
$dCell.Style = "Normal"
$ws.Range("E48").Value = '  +3.61%  '
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9983'
$dCell.Style = "Normal"
$ws.Range("E49").Value = '  -0.66%  '
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '141.82'
$dCell.Style = "Normal"
$ws.Range("E50").Value = '  +0.44%  '
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '0.08359'
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  +3.48%  '
